$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 614, shifting existing rows 614:655 down to 615:656
$ws.Rows.Item(614).Insert()

# Force column A to be treated as plain text so the date-like string is not
# auto-converted into a date serial number (matches the other rows, which
# store the date as a literal text string rather than a real date value).
$ws.Range("A614").NumberFormat = "@"
$ws.Range("A614").Value = "2026/01/09"
$ws.Range("A614").Style = "Normal"

$ws.Range("B614").Value = "金"
$ws.Range("C614").Value = 20
$ws.Range("D614").Value = 201
